# Apply the "Add files via upload" edit:
#  - November sheet: fill in the Column F (second time-slot) attendance figures
#    that mirror the existing Column E numbers, update the 09:00-12:00 time
#    label (shared by D11/E11/F11), and move the active-cell selection.
#  - Overall Attendance sheet: zoom back out to 85%.
# The "Overall Attendance" sheet's numbers are driven entirely by formulas
# referencing November!J.. so they recalculate automatically once the
# November sheet values change.

$wb  = $excel.ActiveWorkbook
$nov = $wb.Worksheets.Item("November")
$oa  = $wb.Worksheets.Item("Overall Attendance")

# --- November: row 10 (Date row) - new class on 20th ---
$nov.Range("F10").Value = 20

# --- November: row 11 (Time row) - new class time label, also applied to
#     the existing D11/E11 cells since they share the same text ---
$newTime = "09:00 TO`n12:00"
$nov.Range("D11").Value = $newTime
$nov.Range("E11").Value = $newTime
$nov.Range("F11").Value = $newTime

# --- November: row 12 (Classes row) ---
$nov.Range("F12").Value = 3

# --- November: attendance rows 14-26, column F ---
$nov.Range("F14").Value = 3
$nov.Range("F15").Value = 3
$nov.Range("F16").Value = 3
$nov.Range("F17").Value = 3
$nov.Range("F18").Value = 3
$nov.Range("F19").Value = 3
$nov.Range("F20").Value = 3
$nov.Range("F21").Value = 0
$nov.Range("F22").Value = 3
$nov.Range("F23").Value = 3
$nov.Range("F24").Value = 3
$nov.Range("F25").Value = 3
$nov.Range("F26").Value = 0

# --- Overall Attendance: zoom out to 85% ---
$null = $oa.Activate()
$excel.ActiveWindow.Zoom = 85

# --- November: restore active sheet/tab and move the selection to E16 ---
$null = $nov.Activate()
$null = $nov.Range("E16").Select()
